$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (A1:D1 containing "a column"/"b column"/"Number"/"Is in DB?").
# This shifts all data rows up by one, turning old row 2 into new row 1, etc.
$ws.Rows(1).Delete()

# Correct one of the numbers that was mis-recorded (old row 4 -> new row 3).
$ws.Range("C3").Value = 1111111111

# Add a text (string) mirror of column C into column D for every populated row,
# built via a TEXT() formula that is then converted to a static value so the
# cells end up holding plain text rather than a live formula.
$dataRows = 1,2,3,4,6,7,8,9
foreach ($r in $dataRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Formula = "=TEXT(C$r,""0"")"
}

$rngTop = $ws.Range("D1:D4")
$rngTop.Copy()
$rngTop.PasteSpecial(-4163)

$rngBottom = $ws.Range("D6:D9")
$rngBottom.Copy()
$rngBottom.PasteSpecial(-4163)

$excel.CutCopyMode = $false
$ws.Range("D1:D9").Select()
